$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schema")
$ws.Activate()

# Fix typo'd item codes in the "Items Used" column (col C) of the Schema sheet.
$ws.Range("C3").Value = "E025, E026, E027, E028,E029"
$ws.Range("C18").Value = "A124_05, A124_06, A124_10"
$ws.Range("C19").Value = "A124_03, A124_08, A124_09"

# Move the saved selection from C28 to C4.
$ws.Range("C4").Select()
